$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 42, pushing the existing rows 42-53 down to 43-54.
$ws.Rows("42:42").Insert()

# Populate the newly inserted row 42. Columns A,B,C,E,F,G,H,I,J,K,L,T keep the
# same "template" values as the record that used to sit at row 42 (now at 43);
# D, M, N, O, P, Q, R, S carry the new weekly figures.
$ws.Range("A42").Value = 9
$ws.Range("B42").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C42").Value = "Metropolitana"
$ws.Range("D42").Value = 44476
$ws.Range("E42").Value = 13
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100102
$ws.Range("H42").Value = "Cítricos"
$ws.Range("I42").Value = 100102006
$ws.Range("J42").Value = "Pomelo"
$ws.Range("K42").Value = "Start Ruby"
$ws.Range("L42").Value = "Primera"
$ws.Range("M42").Value = 350
$ws.Range("N42").Value = 9000
$ws.Range("O42").Value = 9000
$ws.Range("P42").Value = 9000
$ws.Range("Q42").Value = "$/caja 14 kilos empedrada"
$ws.Range("R42").Value = "Provincia de Limarí"
$ws.Range("S42").Value = 643
$ws.Range("T42").Value = 14
